# Update the "K" column (G) values in the save_data sheet.
# These are recalculated strike-count values (K) replacing the old Strike# values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    4  = 3
    5  = 2
    6  = 1
    7  = 4
    8  = 3
    9  = 2
    10 = 2
    11 = 4
    12 = 6
    13 = 6
    14 = 6
    15 = 3
    16 = 4
    17 = 7
    18 = 6
    19 = 7
    20 = 8
    21 = 3
    22 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
